# geração de análises seriais
# Re-generation of the ranking analyses shuffles how tied UFs land on
# particular rows. Update the two sheets whose displayed rankings change:
# "qtd" (count ranking) and "tx-sucesso" (success-rate ranking).

$wb = $excel.ActiveWorkbook

# --- Sheet "qtd": ties at qtd=30 (BA/PB), qtd=4 (MA/RO/AM), qtd=2 (SE/MT) ---
$wsQtd = $wb.Worksheets.Item("qtd")
$wsQtd.Range("A10").Value = "PB"
$wsQtd.Range("A11").Value = "BA"
$wsQtd.Range("A20").Value = "RO"
$wsQtd.Range("A21").Value = "AM"
$wsQtd.Range("A22").Value = "MA"
$wsQtd.Range("A23").Value = "MT"
$wsQtd.Range("A24").Value = "SE"

# --- Sheet "tx-sucesso": ties at txsucesso=100 (MT/RO/MA/XX) ---
$wsTx = $wb.Worksheets.Item("tx-sucesso")
$wsTx.Range("A3").Value = "MA"
$wsTx.Range("A4").Value = "XX"
$wsTx.Range("A5").Value = "RO"
